$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H17").Value = "x"
$ws.Range("I17").Select()
